# Changes From 05 May.xlsx - "Added entries. - Smitha"
#
# Adds new admin/workflow-import related entries to the tracking sheet and
# updates the view/selection to match the author's last edit position.
#
# NOTE: the order cells are written in below matters - it determines which
# shared-string slot each new piece of text lands in, and is chosen to
# mirror the original edit (I9/I10 first, reclaiming the now-unused old
# slot 28 text, then I5, then G9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 / Row 10 text ---------------------------------------------------
# I9: brand-new cell - short note about the admin-driven intake script.
$ws.Range("I9").Value = "Import the intake script from the admin."

# I10 previously held "Import the workflow from admin and release the
# process."; that sentence is now expanded with the "Or / importworkflows"
# command-line alternative, and the cell becomes wrap-text formatted.
$i10 = @'
Import the workflow from admin and release the process.
Or 
importworkflows -Dworkflow.dir=  -Doverwrite=true
'@
$ws.Range("I10").Value = $i10
$ws.Range("I10").WrapText = $true

# --- Row 5 ------------------------------------------------------------
# I5 was a blank (but wrap-styled) cell; it now documents the "admin" route.
$ws.Range("I5").Value = "Import using the admin"

# --- Row 9 continued --------------------------------------------------
# G9: brand-new cell holding the SQL used to set up the Application
# Follow-up work queue / milestone configuration.
$sql = @'
INSERT INTO WORKQUEUE (ADMINISTRATORUSERNAME, ALLOWUSERSUBSCRIPTIONIND, COMMENTS, LASTWRITTEN, NAME, SENSITIVITY, UPPERNAME, VERSIONNO, WORKQUEUEID) VALUES ('admin', '1', 'This work queue is used to assign tasks to case workers when the applications have exceeded 90 days.', '2001-01-01 00:00:00', 'Application Follow-up: Applications that have exceeded 90 days', '1', 'APPLICATION FOLLOW-UP: APPLICATIONS THAT HAVE EXCEEDED 90 DAYS', 1, 45012);
INSERT INTO ALLOCATIONTARGETITEM (ALLOCATIONTARGETID, ALLOCATIONTARGETITEMID, RELATEDID, RELATEDNAME, TYPE) VALUES ('ApplicationFollowUp', 45005, 45012, 'ApplicationFollowUp', 'RL23');
INSERT INTO ALLOCATIONTARGET (ALLOCATIONTARGETID, COMMENTS, NAME) VALUES ('ApplicationFollowUp', 'Application follow up work queue for the case worker.', 'ApplicationFollowUp');
update milestoneconfiguration set duration=90 where milestoneConfigurationID=45001;
INSERT INTO WORKQUEUE (ADMINISTRATORUSERNAME, ALLOWUSERSUBSCRIPTIONIND, COMMENTS, LASTWRITTEN, NAME, SENSITIVITY, UPPERNAME, VERSIONNO, WORKQUEUEID) VALUES ('admin', '1', 'This work queue is used to assign tasks to case workers when the applications have exceeded 90 days.', '2001-01-01 00:00:00', 'متابعة الطلب : الطلبات التي تجاوزت 90 يوما', '1', 'APPLICATION FOLLOW-UP: APPLICATIONS THAT HAVE EXCEEDED 90 DAYS', 1, 45012);
'@
$ws.Range("G9").Value = $sql
$ws.Range("G9").WrapText = $true

# --- View / selection -----------------------------------------------------
$ws.Activate()
$ws.Range("D1").Select()
$ws.Range("G10").Select()
